$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.891.97'
$ws.Range("E2").Value = '  +2.03%  '
$ws.Range("D3").Value = '3.048.81'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''525.76'
$ws.Range("D6").Value = '''142.70'
$ws.Range("E6").Value = '  +5.98%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.445'
$ws.Range("E8").Value = '  +4.76%  '
$ws.Range("E9").Value = '  +4.03%  '
$ws.Range("E10").Value = '  +8.98%  '
$ws.Range("D11").Value = '''0.370'
$ws.Range("E11").Value = '  +5.76%  '
$ws.Range("E12").Value = '  +2.67%  '
$ws.Range("D13").Value = '3.569.76'
$ws.Range("E13").Value = '  +2.23%  '
$ws.Range("D14").Value = '''26.87'
$ws.Range("E14").Value = '  +8.73%  '
$ws.Range("E15").Value = '  +17.40%  '
$ws.Range("D16").Value = '57.816.37'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '''6.22'
$ws.Range("E17").Value = '  +7.46%  '
$ws.Range("D18").Value = '3.046.79'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("D19").Value = '''12.93'
$ws.Range("E20").Value = '  +6.16%  '
$ws.Range("D21").Value = '''342.64'
$ws.Range("E21").Value = '  +5.67%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '''0.500'
$ws.Range("E23").Value = '  +7.86%  '
$ws.Range("D24").Value = '''65.31'
$ws.Range("E24").Value = '  +6.56%  '
$ws.Range("D25").Value = '''0.173'
$ws.Range("E25").Value = '  +6.93%  '
$ws.Range("D26").Value = '0.0₃0971'
$ws.Range("E26").Value = '  +7.39%  '
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").Value = '''7.02'
$ws.Range("E28").Value = '  +8.21%  '
$ws.Range("D29").Value = '''7.30'
$ws.Range("E29").Value = '  +7.76%  '
$ws.Range("E30").Value = '  +7.61%  '
$ws.Range("E31").Value = '  +5.84%  '
$ws.Range("D32").Value = '''21.08'
$ws.Range("E32").Value = '  +5.97%  '
$ws.Range("D33").Value = '''156.50'
$ws.Range("E33").Value = '  +1.29%  '
$ws.Range("D34").Value = '''4.74'
$ws.Range("E34").Value = '  +6.26%  '
$ws.Range("E35").Value = '  +5.74%  '
$ws.Range("D37").Value = '''25.82'
$ws.Range("E37").Value = '  +10.08%  '
$ws.Range("D38").Value = '''0.0697'
$ws.Range("E38").Value = '  +3.79%  '
$ws.Range("D39").Value = '3.080.44'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").Value = '''37.71'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("E41").Value = '  +8.59%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  +5.15%  '
$ws.Range("D44").Value = '''0.662'
$ws.Range("E44").Value = '  +4.04%  '
$ws.Range("D45").Value = '2.324.42'
$ws.Range("E45").Value = '  +5.54%  '
$ws.Range("E46").Value = '  +1.64%  '
$ws.Range("E47").Value = '  +4.81%  '
$ws.Range("E48").Value = '  +5.05%  '
$ws.Range("E49").Value = '  +5.90%  '
$ws.Range("D50").Value = '''20.11'
$ws.Range("E50").Value = '  +5.20%  '
$ws.Range("D51").Value = '''0.0895'
$ws.Range("E51").Value = '  +6.07%  '
